$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Supplier" column header in K1
$ws.Range("K1").Value2 = "Supplier"

# Fill the Supplier column (id 2) for every data row, matching the
# default/global supplier used by the rest of the import fixture
$ws.Range("K2:K7").Value2 = 2

# Row 2 content grew (new column), row autofit shrank its height slightly
$ws.Rows.Item(2).RowHeight = 15

# Leave the selection on the newly added K7 cell, scrolled so column E is
# the first visible column (mirrors the saved view state of the edited file)
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K7").Select() | Out-Null
